$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the tab-ratio on the book's window (workbookView tabRatio 986 -> 989).
$win = $wb.Windows.Item(1)
$win.TabRatio = 0.989

# Add the opening/closing time labels + placeholders below the existing
# employee-list template rows.
$ws.Range("A7").Value = '始業時間'
$ws.Range("B7").Value = '${openingTime}'
$ws.Range("A8").Value = '終業時間'
$ws.Range("B8").Value = '${closingTime}'

# Drop the now-unneeded explicit widths on column A and columns AG onward,
# leaving only the narrow day-of-month columns B:AF customised.
$ws.Columns.Item(1).ClearFormats()
$ws.Range($ws.Cells.Item(1, 33), $ws.Cells.Item(1, 1025)).EntireColumn.ClearFormats()
$ws.Columns("B:AF").ColumnWidth = 2.3030303030303054

# Reset the saved selection back to A1.
$ws.Range("A1").Select()
